$wb = $excel.ActiveWorkbook

# ALC row 62  (@@ -3706,25 +3706,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 29419476
$ws.Range("I62").Value = 7864
$ws.Range("J62").Value = 166673660
$ws.Range("K62").Value = 7864
$ws.Range("L62").Value = 166673660
$ws.Range("M62").Value = -7240
$ws.Range("N62").Value = -166674908

# ALC row 65  (@@ -3856,25 +3856,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H65").Value = 29419476
$ws.Range("I65").Value = 7864
$ws.Range("J65").Value = 166673660
$ws.Range("K65").Value = 39320
$ws.Range("L65").Value = 833368300
$ws.Range("M65").Value = -36200
$ws.Range("N65").Value = -833374540

# ALC row 98  (@@ -5500,25 +5500,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 297875
$ws.Range("I98").Value = 1826.7037
$ws.Range("J98").Value = 1439775.6
$ws.Range("K98").Value = 1826.7037
$ws.Range("L98").Value = 1439775.6
$ws.Range("M98").Value = -328.7037
$ws.Range("N98").Value = -1442771.6

# ALC row 122  (@@ -6688,25 +6688,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 297875
$ws.Range("I122").Value = 1826.7037
$ws.Range("J122").Value = 1439775.6
$ws.Range("K122").Value = 5480.1111
$ws.Range("L122").Value = 4319326.800000001
$ws.Range("M122").Value = -3030.1111
$ws.Range("N122").Value = -4324226.800000001

# ALC row 125  (@@ -6835,22 +6835,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 2687.5715
$ws.Range("I125").Value = 2668.75
$ws.Range("K125").Value = 24018.75
$ws.Range("M125").Value = -21558.75

# ALC row 131  (@@ -7129,22 +7129,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 27720.889
$ws.Range("I131").Value = 1859.9615
$ws.Range("K131").Value = 5579.8845
$ws.Range("M131").Value = -539.8845000000001

# ALC row 132  (@@ -7181,22 +7181,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1341.8478
$ws.Range("I132").Value = 1241.7949
$ws.Range("K132").Value = 3725.384700000001
$ws.Range("M132").Value = -1195.384700000001

# ALC row 137  (@@ -7429,25 +7429,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 10930.944
$ws.Range("I137").Value = 28125
$ws.Range("J137").Value = 6018.357
$ws.Range("K137").Value = 84375
$ws.Range("L137").Value = 18055.071
$ws.Range("M137").Value = -81825
$ws.Range("N137").Value = -23155.071

# ALC row 138  (@@ -7481,25 +7481,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 1792.21
$ws.Range("I138").Value = 1236.9302
$ws.Range("J138").Value = 2211.1052
$ws.Range("K138").Value = 3710.7906
$ws.Range("L138").Value = 6633.3156
$ws.Range("M138").Value = 1429.2094
$ws.Range("N138").Value = -16913.3156

# ALC row 141  (@@ -7628,22 +7628,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 7998
$ws.Range("I141").Value = 7998
$ws.Range("K141").Value = 23994
$ws.Range("M141").Value = -18814

# ARM row 61  (@@ -10644,22 +10644,22 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7088.593
$ws.Range("I61").Value = 5896.9473
$ws.Range("K61").Value = 5896.9473
$ws.Range("M61").Value = -5684.9473

# ARM row 88  (@@ -11955,25 +11955,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2212.9092
$ws.Range("I88").Value = 1972.6
$ws.Range("J88").Value = 2413.1667
$ws.Range("K88").Value = 1972.6
$ws.Range("L88").Value = 2413.1667
$ws.Range("M88").Value = -1566.6
$ws.Range("N88").Value = -3225.1667

# ARM row 91  (@@ -12102,25 +12102,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2212.9092
$ws.Range("I91").Value = 1972.6
$ws.Range("J91").Value = 2413.1667
$ws.Range("K91").Value = 1972.6
$ws.Range("L91").Value = 2413.1667
$ws.Range("M91").Value = -568.5999999999999
$ws.Range("N91").Value = -5221.1667

# ARM row 107  (@@ -12880,19 +12880,22 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H107").Value = 199999
$ws.Range("J107").Value = 199999
$ws.Range("L107").Value = 199999
$ws.Range("N107").Value = -207679

# ARM row 122  (@@ -13600,22 +13603,22 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 10331
$ws.Range("I122").Value = 6918.5884
$ws.Range("K122").Value = 20755.7652
$ws.Range("M122").Value = -18305.7652

# ARM row 136  (@@ -14274,22 +14277,22 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 7088.593
$ws.Range("I136").Value = 5896.9473
$ws.Range("K136").Value = 17690.8419
$ws.Range("M136").Value = -15140.8419

# BSM row 20  (@@ -15556,25 +15559,25 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2482.3735
$ws.Range("I20").Value = 1981.5077
$ws.Range("J20").Value = 3734.5386
$ws.Range("K20").Value = 1981.5077
$ws.Range("L20").Value = 3734.5386
$ws.Range("M20").Value = -1734.5077
$ws.Range("N20").Value = -4228.5386

# BSM row 99  (@@ -19415,25 +19418,25 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 7738.207
$ws.Range("I99").Value = 8722.9
$ws.Range("J99").Value = 5550
$ws.Range("K99").Value = 8722.9
$ws.Range("L99").Value = 5550
$ws.Range("M99").Value = -7224.9
$ws.Range("N99").Value = -8546

# CRP row 122  (@@ -27442,22 +27445,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2176.2727
$ws.Range("I122").Value = 1571.875
$ws.Range("K122").Value = 4715.625
$ws.Range("M122").Value = -2265.625

# CRP row 132  (@@ -27929,22 +27932,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3474.92
$ws.Range("I132").Value = 3242.682
$ws.Range("K132").Value = 9728.045999999998
$ws.Range("M132").Value = -7198.045999999998

# CRP row 134  (@@ -28030,25 +28033,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 12513.5
$ws.Range("I134").Value = 15090.125
$ws.Range("J134").Value = 2207
$ws.Range("K134").Value = 45270.375
$ws.Range("L134").Value = 6621
$ws.Range("M134").Value = -42735.375
$ws.Range("N134").Value = -11691

# CUL row 4  (@@ -28620,22 +28623,22 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4113672.8
$ws.Range("I4").Value = 4113672.8
$ws.Range("K4").Value = 12341018.4
$ws.Range("M4").Value = -12340906.4

# CUL row 112  (@@ -33948,19 +33951,22 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 3526
$ws.Range("I112").Value = 3526
$ws.Range("K112").Value = 10578
$ws.Range("M112").Value = -9470

# CUL row 113  (@@ -33994,25 +34000,25 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 2668
$ws.Range("I113").Value = 2550
$ws.Range("J113").Value = 2684.8572
$ws.Range("K113").Value = 7650
$ws.Range("L113").Value = 8054.571599999999
$ws.Range("M113").Value = -5480
$ws.Range("N113").Value = -12394.5716

# CUL row 137  (@@ -35209,22 +35215,22 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 22224364
$ws.Range("I137").Value = 1248.4286
$ws.Range("K137").Value = 3745.2858
$ws.Range("M137").Value = 1354.7142

# CUL row 140  (@@ -35365,25 +35371,25 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 6580565
$ws.Range("I140").Value = 22728020
$ws.Range("J140").Value = 1972.5186
$ws.Range("K140").Value = 68184060
$ws.Range("L140").Value = 5917.5558
$ws.Range("M140").Value = -68178880
$ws.Range("N140").Value = -16277.5558

# GSM row 5  (@@ -35716,22 +35722,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1799.5
$ws.Range("I5").Value = 1799.5
$ws.Range("K5").Value = 1799.5
$ws.Range("M5").Value = -1687.5

# GSM row 70  (@@ -38883,25 +38889,25 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 16748.5
$ws.Range("J70").Value = 16597.8
$ws.Range("L70").Value = 16597.8
$ws.Range("N70").Value = -17137.8

# GSM row 73  (@@ -39030,25 +39036,25 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 16748.5
$ws.Range("J73").Value = 16597.8
$ws.Range("L73").Value = 16597.8
$ws.Range("N73").Value = -18469.8

# GSM row 96  (@@ -40148,22 +40154,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 54199.4
$ws.Range("J96").Value = 54199.4
$ws.Range("L96").Value = 54199.4
$ws.Range("N96").Value = -59691.4

# GSM row 122  (@@ -41410,22 +41416,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2599.625
$ws.Range("I122").Value = 1489.8
$ws.Range("K122").Value = 4469.4
$ws.Range("M122").Value = -2019.4

# GSM row 126  (@@ -41603,25 +41609,25 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 10632.4
$ws.Range("J126").Value = 10115.818
$ws.Range("L126").Value = 30347.454
$ws.Range("N126").Value = -35287.454

# GSM row 132  (@@ -41888,22 +41894,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 44578.36
$ws.Range("I132").Value = 47925.41
$ws.Range("K132").Value = 143776.23
$ws.Range("M132").Value = -141246.23

# GSM row 135  (@@ -42038,22 +42044,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 68200
$ws.Range("J135").Value = 68200
$ws.Range("L135").Value = 68200
$ws.Range("N135").Value = -78340

# LTW row 7  (@@ -42711,25 +42717,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7870.1333
$ws.Range("J7").Value = 7699.75
$ws.Range("L7").Value = 7699.75
$ws.Range("N7").Value = -7923.75

# LTW row 126  (@@ -48467,25 +48473,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 7870.1333
$ws.Range("J126").Value = 7699.75
$ws.Range("L126").Value = 23099.25
$ws.Range("N126").Value = -28039.25

# LTW row 132  (@@ -48758,22 +48764,22 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 8055.2285
$ws.Range("I132").Value = 4141.077
$ws.Range("K132").Value = 12423.231
$ws.Range("M132").Value = -9893.231

# LTW row 133  (@@ -48810,22 +48816,22 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H133").Value = 69999
$ws.Range("J133").Value = 69999
$ws.Range("L133").Value = 69999
$ws.Range("N133").Value = -75059

# LTW row 136  (@@ -48957,22 +48963,22 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2876.0527
$ws.Range("I136").Value = 2679.9
$ws.Range("K136").Value = 8039.700000000001
$ws.Range("M136").Value = -5489.700000000001

# WVR row 2  (@@ -49339,22 +49345,22 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 265919.5
$ws.Range("I2").Value = 265919.5
$ws.Range("K2").Value = 265919.5
$ws.Range("M2").Value = -265807.5

# WVR row 62  (@@ -52258,25 +52264,25 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8310.105
$ws.Range("I62").Value = 5993.7
$ws.Range("J62").Value = 10883.889
$ws.Range("K62").Value = 5993.7
$ws.Range("L62").Value = 10883.889
$ws.Range("M62").Value = -5369.7
$ws.Range("N62").Value = -12131.889

# WVR row 65  (@@ -52408,25 +52414,25 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 8310.105
$ws.Range("I65").Value = 5993.7
$ws.Range("J65").Value = 10883.889
$ws.Range("K65").Value = 29968.5
$ws.Range("L65").Value = 54419.44499999999
$ws.Range("M65").Value = -26848.5
$ws.Range("N65").Value = -60659.44499999999

# WVR row 126  (@@ -55367,25 +55373,25 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 4903.1113
$ws.Range("I126").Value = 2748.75
$ws.Range("J126").Value = 6626.6
$ws.Range("K126").Value = 8246.25
$ws.Range("L126").Value = 19879.8
$ws.Range("M126").Value = -5776.25
$ws.Range("N126").Value = -24819.8

# WVR row 132  (@@ -55661,25 +55667,22 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2798.6667
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2798.6667
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 8396.000100000001
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -13456.0001
